# The post "「気をつけて　ジュード」" (previously row 363) was removed from the
# spreadsheet. Deleting the entire worksheet row shifts every following row
# up by one, which matches the diff (row 364 -> 363, ..., row 484 -> 483,
# and the sheet's used range shrinks from A1:C484 to A1:C483).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(363).Delete()
